$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref="D2"; Val="26.156.18"},
    @{Ref="E2"; Val="  -4.51%  "},
    @{Ref="D3"; Val="1.656.45"},
    @{Ref="E3"; Val="  -3.31%  "},
    @{Ref="D4"; Val="1.005"},
    @{Ref="E4"; Val="  +0.00%  "},
    @{Ref="D5"; Val="215.52"},
    @{Ref="E5"; Val="  -4.04%  "},
    @{Ref="D6"; Val="0.5110"},
    @{Ref="E6"; Val="  -3.70%  "},
    @{Ref="E7"; Val="  +0.11%  "},
    @{Ref="D8"; Val="0.2581"},
    @{Ref="E8"; Val="  -3.28%  "},
    @{Ref="D9"; Val="0.06420"},
    @{Ref="E9"; Val="  -4.13%  "},
    @{Ref="D10"; Val="19.94"},
    @{Ref="E10"; Val="  -4.59%  "},
    @{Ref="D11"; Val="0.07802"},
    @{Ref="E11"; Val="  +1.49%  "},
    @{Ref="D12"; Val="1.661.04"},
    @{Ref="E12"; Val="  -2.86%  "},
    @{Ref="D13"; Val="4.286"},
    @{Ref="E13"; Val="  -4.98%  "},
    @{Ref="D14"; Val="1.882.74"},
    @{Ref="E14"; Val="  -3.39%  "},
    @{Ref="D15"; Val="0.5518"},
    @{Ref="E15"; Val="  -5.40%  "},
    @{Ref="D16"; Val="0.0₅8023"},
    @{Ref="E16"; Val="  -2.62%  "},
    @{Ref="D17"; Val="64.07"},
    @{Ref="E17"; Val="  -5.97%  "},
    @{Ref="D18"; Val="26.169.23"},
    @{Ref="E18"; Val="  -4.43%  "},
    @{Ref="D20"; Val="209.42"},
    @{Ref="E20"; Val="  -6.06%  "},
    @{Ref="D21"; Val="4.416"},
    @{Ref="E21"; Val="  -4.71%  "},
    @{Ref="D22"; Val="10.07"},
    @{Ref="E22"; Val="  -3.13%  "},
    @{Ref="D23"; Val="6.016"},
    @{Ref="E23"; Val="  +0.04%  "},
    @{Ref="E24"; Val="  +0.14%  "},
    @{Ref="D25"; Val="143.64"},
    @{Ref="E25"; Val="  -0.68%  "},
    @{Ref="D26"; Val="1.741"},
    @{Ref="E26"; Val="  +3.06%  "},
    @{Ref="D27"; Val="0.1178"},
    @{Ref="E27"; Val="  -2.68%  "},
    @{Ref="D28"; Val="6.979"},
    @{Ref="E28"; Val="  -3.81%  "},
    @{Ref="D29"; Val="15.82"},
    @{Ref="E29"; Val="  -2.84%  "},
    @{Ref="D30"; Val="0.05126"},
    @{Ref="E30"; Val="  -4.71%  "},
    @{Ref="E31"; Val="  -4.03%  "},
    @{Ref="D32"; Val="3.352"},
    @{Ref="E32"; Val="  -3.51%  "},
    @{Ref="D33"; Val="3.222"},
    @{Ref="E33"; Val="  -6.45%  "},
    @{Ref="D34"; Val="1.566"},
    @{Ref="E34"; Val="  -4.65%  "},
    @{Ref="D35"; Val="2.749"},
    @{Ref="E35"; Val="  -4.15%  "},
    @{Ref="E36"; Val="  -0.82%  "},
    @{Ref="D37"; Val="0.9286"},
    @{Ref="E37"; Val="  -2.49%  "},
    @{Ref="D38"; Val="0.5686"},
    @{Ref="E38"; Val="  -2.85%  "},
    @{Ref="D39"; Val="1.159.78"},
    @{Ref="E39"; Val="  +5.97%  "},
    @{Ref="D40"; Val="0.01590"},
    @{Ref="E40"; Val="  -2.88%  "},
    @{Ref="D41"; Val="2.556"},
    @{Ref="E41"; Val="  -0.06%  "},
    @{Ref="D42"; Val="1.005"},
    @{Ref="D43"; Val="0.8342"},
    @{Ref="E43"; Val="  -0.96%  "},
    @{Ref="D44"; Val="5.635"},
    @{Ref="E44"; Val="  -2.82%  "},
    @{Ref="D45"; Val="100.44"},
    @{Ref="E45"; Val="  -0.35%  "},
    @{Ref="D46"; Val="1.792.97"},
    @{Ref="E46"; Val="  -3.37%  "},
    @{Ref="D47"; Val="0.0₈117"},
    @{Ref="E47"; Val="  +0.08%  "},
    @{Ref="D48"; Val="0.4547"},
    @{Ref="E48"; Val="  +0.21%  "},
    @{Ref="D49"; Val="55.75"},
    @{Ref="E49"; Val="  -3.62%  "},
    @{Ref="D50"; Val="1.004"},
    @{Ref="E50"; Val="  -0.02%  "},
    @{Ref="D51"; Val="7.867"},
    @{Ref="E51"; Val="  -2.66%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Val
    $cell.NumberFormat = "General"
}
